$wb = $excel.ActiveWorkbook

$wsSize = $wb.Worksheets.Item("Size")
$wsSize.Range("G2").Value = 0
$wsSize.Range("G3").Value = 0
$wsSize.Range("G4").Value = 169.892
$wsSize.Range("C5").Value = 262.3528836842104
$wsSize.Range("D5").Value = 1157.10682
$wsSize.Range("E5").Value = 40.19118066666667
$wsSize.Range("F5").Value = 0
$wsSize.Range("G5").Value = 1459.650884350877

$wsCost = $wb.Worksheets.Item("Cost")
$wsCost.Range("I2").Value = 0.2878138498753487
$wsCost.Range("I3").Value = 0
$wsCost.Range("I4").Value = 0
$wsCost.Range("I5").Value = 0.0339784
$wsCost.Range("E6").Value = 0.02623528836842104
$wsCost.Range("F6").Value = 0.115710682
$wsCost.Range("G6").Value = 0.004019118066666667
$wsCost.Range("H6").Value = 0
$wsCost.Range("I6").Value = 0.1459650884350877
$wsCost.Range("I7").Value = 0
$wsCost.Range("I8").Value = 0
$wsCost.Range("I9").Value = 0.04742910670265638
$wsCost.Range("E10").Value = 0.005493120444462371
$wsCost.Range("F10").Value = 0.02422739571263718
$wsCost.Range("G10").Value = 0.0008415192282501857
$wsCost.Range("H10").Value = 0
$wsCost.Range("I10").Value = 0.03056203538534974
$wsCost.Range("I11").Value = 0.01018746213436753
$wsCost.Range("E12").Value = 0.004912450177583252
$wsCost.Range("F12").Value = 0.01456426425675683
$wsCost.Range("G12").Value = 0.0002150427835472044
$wsCost.Range("H12").Value = 0
$wsCost.Range("I12").Value = 0.01969175721788729

$wsInd = $wb.Worksheets.Item("Indicators")
$wsInd.Range("C2").Value = 430.8938900000001
$wsInd.Range("D2").Value = 427.4078603793104
$wsInd.Range("E2").Value = 858.3017503793104
$wsInd.Range("E3").Value = 1
$wsInd.Range("E4").Value = 0
$wsInd.Range("C5").Value = 0.0002999999999999999
$wsInd.Range("D5").Value = 0.0008787878787878787
$wsInd.Range("E5").Value = 0.0005882185534130422
$wsInd.Range("E6").Value = 2.450429441131381
